$d = $word.ActiveDocument

# --- Edit 1: merge the "@EN" + "D" runs (HandmadePageBreak section) into one "@END" run ---
# The literal text "@END" appears many times in this document (once per lint
# rule block), so scope the Find/Replace to just the paragraph that follows
# the (unique) "Перевод на новую страницу..." explanation text.
$anchor = $d.Content
$anchor.Find.Execute("Перевод на новую страницу") | Out-Null
$anchorPara = $anchor.Paragraphs(1)
$splitPara = $anchorPara.Next()
$splitRange = $d.Range($splitPara.Range.Start, $splitPara.Range.End)
$splitRange.Find.Execute("@END", $true, $false, $false, $false, $false, `
    $true, 1, $false, "@END", 2) | Out-Null

# --- Edit 2: append the new "IncorrectHeadingText" lint block at the end ---
# The document currently ends with a paragraph containing "@END" (closing the
# IncorrectFontSize block). Insert four new paragraphs after it:
#   @END
#   (empty paragraph)
#   @BEGIN IncorrectHeadingText
#   Правильное оформление заголовка — «$Expected».
$last = $d.Paragraphs.Last
$last.Range.InsertParagraphAfter()
$p1 = $d.Paragraphs.Last
$p1.Range.Text = "@END"

$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Last

$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Last
$p3.Range.Text = "@BEGIN IncorrectHeadingText"

$p3.Range.InsertParagraphAfter()
$p4 = $d.Paragraphs.Last
$p4.Range.Text = "Правильное оформление заголовка "
$p4.Range.InsertAfter("— «`$Expected».")

$p4full = $d.Paragraphs.Last.Range
$tailFind = $d.Range($p4full.Start, $p4full.End)
$tailFind.Find.Execute("— «") | Out-Null
Write-Host "tailFind after find, start:" $tailFind.Start "end:" $tailFind.End
$tail = $d.Range($tailFind.Start, $p4full.End - 1)
Write-Host "tail text:" $tail.Text
$tail.LanguageID = "ru-RU"
